# Apply updated odds values to Sheet1 as per the 2025-02-05 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.03
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.93
$ws.Range("T2").Value = 1.3
$ws.Range("W2").Value = 1.8
$ws.Range("X2").Value = 1.91

# Row 3
$ws.Range("G3").Value = 5.75
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 1.6
$ws.Range("M3").Value = 1.03
$ws.Range("O3").Value = 1.22
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 3
$ws.Range("W3").Value = 1.83
$ws.Range("X3").Value = 1.83
$ws.Range("AF3").Value = 7

# Row 4
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.9
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.1
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.03
$ws.Range("O4").Value = 1.19
$ws.Range("T4").Value = 1.37
$ws.Range("Z4").Value = 8
$ws.Range("AC4").Value = 12
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 17

# Row 5
$ws.Range("J5").Value = 1.24
$ws.Range("K5").Value = 3.6
$ws.Range("Q5").Value = 1.2
$ws.Range("R5").Value = 4.2
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.12
$ws.Range("Y5").Value = 10.25
$ws.Range("Z5").Value = 6.2
$ws.Range("AA5").Value = 13
$ws.Range("AB5").Value = 5.2
$ws.Range("AC5").Value = 10.75
$ws.Range("AD5").Value = 40
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 55
$ws.Range("AH5").Value = 250
$ws.Range("AI5").Value = 120
$ws.Range("AN5").Value = 450

# Row 7
$ws.Range("G7").Value = 1.11
$ws.Range("H7").Value = 9
$ws.Range("I7").Value = 21
$ws.Range("K7").Value = 3.6
$ws.Range("W7").Value = 1.95
$ws.Range("X7").Value = 1.8
$ws.Range("AJ7").Value = 101
$ws.Range("AM7").Value = 126
$ws.Range("AN7").Value = 81

# Row 8
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 3.2
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 15
$ws.Range("Q8").Value = 1.67
$ws.Range("R8").Value = 2.15
$ws.Range("Z8").Value = 13
$ws.Range("AA8").Value = 10
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 19
$ws.Range("AD8").Value = 23
$ws.Range("AJ8").Value = 15
$ws.Range("AL8").Value = 26
$ws.Range("AM8").Value = 21
$ws.Range("AN8").Value = 26

# Row 9
$ws.Range("G9").Value = 1.33
$ws.Range("H9").Value = 5.5
$ws.Range("I9").Value = 7.5
$ws.Range("J9").Value = 1.8
$ws.Range("K9").Value = 2.6
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 17
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.57
$ws.Range("R9").Value = 2.35
$ws.Range("S9").Value = 2.38
$ws.Range("T9").Value = 1.53
$ws.Range("U9").Value = 1.29
$ws.Range("V9").Value = 3.5
$ws.Range("W9").Value = 1.83
$ws.Range("X9").Value = 1.83
$ws.Range("Y9").Value = 8
$ws.Range("Z9").Value = 7
$ws.Range("AB9").Value = 8.5
$ws.Range("AC9").Value = 11
$ws.Range("AD9").Value = 26
$ws.Range("AE9").Value = 17
$ws.Range("AF9").Value = 10
$ws.Range("AI9").Value = 21
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 81
$ws.Range("AN9").Value = 51
$ws.Range("AO9").Value = 301
$ws.Range("AP9").Value = 1.98
$ws.Range("AQ9").Value = 1.88
